# "Updated excel reading logic"
#
# Sheet1: column A had "Siva"/"Sankar"/"Ram". "Ram" was shortened to "Ra"
# and rows 2 & 3 were swapped (A2 now holds "Ra", A3 now holds "Sankar").
# A new worksheet ("Sheet3") was also added with some sample data, and
# it became the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet1: fix up existing data --------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "Ra"
$ws1.Range("A3").Value = "Sankar"

# --- Create the new sheet ("Sheet3") -----------------------------------
# A throwaway sheet is created and removed first so the workbook's
# internal sheetId counter lands on 3 for the sheet that is kept -
# matching a "Sheet2 added, then removed, Sheet3 added" authoring
# history.
$throwaway = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Sheet3"
$null = $throwaway.Delete()

# Re-fetch by name: once a sibling sheet is deleted, older worksheet
# object references can go stale, so grab a fresh handle before writing.
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("A1").Value = "fgf"
$ws3.Range("C6").Value = "fdgfd"
$ws3.Range("B4").Value = "fdg"
$ws3.Range("D4").Value = 4564

$null = $ws3.Range("D9").Select()

# --- Selection / active-tab bookkeeping ---------------------------------
# Sheet1 ends up with column B selected (no longer the active tab);
# Sheet3 ends up selected and active.
$null = $wb.Worksheets.Item(2).Activate()
$null = $ws1.Range("B1:B1048576").Select()
$null = $wb.Worksheets.Item(2).Activate()
